# Refresh the cryptocurrency price ("Price", column D) and 1-hour
# volume-change ("Volume(1h)", column E) figures for each coin row,
# as produced by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''71.066.48'
$ws.Range("E2").Value = '  -1.58%  '

$ws.Range("D3").Value = '''2.569.29'
$ws.Range("E3").Value = '  -5.27%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = '''582.02'
$ws.Range("E5").Value = '  -2.98%  '

$ws.Range("D6").Value = '''170.97'
$ws.Range("E6").Value = '  -2.96%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '''0.512'
$ws.Range("E8").Value = '  -2.49%  '

$ws.Range("D9").Value = '''2.568.20'
$ws.Range("E9").Value = '  -5.27%  '

$ws.Range("D10").Value = '''0.167'
$ws.Range("E10").Value = '  -1.47%  '

$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("E12").Value = '  -0.85%  '

$ws.Range("E13").Value = '  -3.24%  '

$ws.Range("D14").Value = '''3.054.98'
$ws.Range("E14").Value = '  -4.79%  '

$ws.Range("E15").Value = '  -0.72%  '

$ws.Range("D16").Value = '''70.903.76'
$ws.Range("E16").Value = '  -1.50%  '

$ws.Range("D17").Value = '''25.26'
$ws.Range("E17").Value = '  -4.16%  '

$ws.Range("D18").Value = '''2.598.41'
$ws.Range("E18").Value = '  -4.17%  '

$ws.Range("D19").Value = '''11.88'
$ws.Range("E19").Value = '  -3.24%  '

$ws.Range("E20").Value = '  -5.17%  '

$ws.Range("D21").Value = '''365.08'
$ws.Range("E21").Value = '  -2.15%  '

$ws.Range("D22").Value = '''4.01'
$ws.Range("E22").Value = '  -3.93%  '

$ws.Range("E23").Value = '  -2.13%  '

$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("D25").Value = '''70.29'
$ws.Range("E25").Value = '  -2.95%  '

$ws.Range("D26").Value = '''4.16'
$ws.Range("E26").Value = '  -4.57%  '

$ws.Range("D27").Value = '''9.26'
$ws.Range("E27").Value = '  -6.29%  '

$ws.Range("D28").Value = '''2.755.25'
$ws.Range("E28").Value = '  -3.38%  '

$ws.Range("E29").Value = '  +0.29%  '

$ws.Range("D30").Value = '''0.0₃0931'
$ws.Range("E30").Value = '  -7.23%  '

$ws.Range("D31").Value = '''7.83'
$ws.Range("E31").Value = '  -3.95%  '

$ws.Range("D32").Value = '''1.32'
$ws.Range("E32").Value = '  +0.92%  '

$ws.Range("D33").Value = '''485.18'
$ws.Range("E33").Value = '  -4.64%  '

$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("D36").Value = '''157.76'
$ws.Range("E36").Value = '  -3.88%  '

$ws.Range("E37").Value = '  +5.03%  '

$ws.Range("D38").Value = '''18.88'
$ws.Range("E38").Value = '  -4.24%  '

$ws.Range("D39").Value = '''18.87'
$ws.Range("E39").Value = '  -1.21%  '

$ws.Range("E40").Value = '  -3.76%  '

$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("E42").Value = '  -5.94%  '

$ws.Range("D43").Value = '''2.49'
$ws.Range("E43").Value = '  -3.39%  '

$ws.Range("E44").Value = '  -5.74%  '

$ws.Range("D45").Value = '''0.320'
$ws.Range("E45").Value = '  -4.28%  '

$ws.Range("D46").Value = '''38.68'
$ws.Range("E46").Value = '  -1.98%  '

$ws.Range("D47").Value = '''147.39'
$ws.Range("E47").Value = '  -6.31%  '

$ws.Range("E48").Value = '  -4.36%  '

$ws.Range("D49").Value = '''0.532'
$ws.Range("E49").Value = '  -5.68%  '

$ws.Range("D50").Value = '''1.64'
$ws.Range("E50").Value = '  -7.55%  '

$ws.Range("D51").Value = '''0.595'
$ws.Range("E51").Value = '  -2.24%  '
